# Physics108_BABYBLUE/connection_diagram.xlsx - "adding data and files from Sat cooldown"
#
# Changes:
#  1. Add a new shared string "squid return meas" (used by B8).
#  2. B8: "squid return" (s=10) -> "squid return meas" (new string).
#  3. B9: "gnd" (s=0) -> "device" (s=1).
#  4. Row 8 grows to a 30pt (wrapped) row height, matching the other
#     header rows (2 and 5) that already use ht="30".
#  5. Active selection on the sheet moves to J4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connection_diagram")

# Row 8: B8 goes from "squid return" to "squid return meas"; row grows to 30pt.
$ws.Range("B8").Value = "squid return meas"
$ws.Rows.Item(8).RowHeight = 30

# Row 9: B9 goes from "gnd" to "device"
$ws.Range("B9").Value = "device"

# Reflect the saved cursor position/selection from the workbook.
$ws.Range("J4").Select()
